# "Change dates in all materials"
#
# The schedule's weekly dates (column A, rows 4-35) were individually
# computed with per-row relative formulas ("=A2+7", "=A3+7", ...). Re-enter
# them as two fill-down ranges so Excel groups them into shared formulas
# (A4:A11 and A12:A35), matching how the author re-did the date column.
# Also bump the row heights for rows 3-4 (taller wrapped text) and make
# sure the re-filled date cells keep the same date/wrap/locked formatting
# already used by the cell just above them (A3), instead of the separate
# (but visually identical) style that used to be applied to A4:A35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row height changes (row 3: 18 -> 35.25, row 4: 32.25 -> 57)
$ws.Rows.Item(3).RowHeight = 35.25
$ws.Rows.Item(4).RowHeight = 57

# Re-enter the weekly "+7 days" formulas a range at a time (fill-down),
# which makes Excel store each block as a single shared formula.
$ws.Range("A4:A11").Formula = "=A2+7"
$ws.Range("A12:A35").Formula = "=A10+7"

# Make the re-filled date cells (A4:A35) pick up the exact same
# number format / wrap / protection already on A3, instead of the
# separate duplicate style they used to reference.
$ws.Range("A3").Copy()
$ws.Range("A4:A35").PasteSpecial(-4122)
$excel.CutCopyMode = $false
